$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.243.43"
$ws.Range("E2").Value = "  -1.47%  "
$ws.Range("D3").Value = "2.296.80"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "533.73"
$ws.Range("E5").Value = "  -3.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.58"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.586"
$ws.Range("E8").Value = "  +2.46%  "
$ws.Range("D9").Value = "2.294.23"
$ws.Range("E9").Value = "  -0.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0997"
$ws.Range("E10").Value = "  -3.00%  "
$ws.Range("E11").Value = "  -1.36%  "
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.330"
$ws.Range("E13").Value = "  -2.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.47"
$ws.Range("E14").Value = "  -1.57%  "
$ws.Range("D15").Value = "2.704.40"
$ws.Range("E15").Value = "  -1.33%  "
$ws.Range("D16").Value = "58.146.96"
$ws.Range("E16").Value = "  -1.67%  "
$ws.Range("E17").Value = "  -1.55%  "
$ws.Range("D18").Value = "2.267.69"
$ws.Range("E18").Value = "  -2.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.55"
$ws.Range("E19").Value = "  -2.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.19"
$ws.Range("E20").Value = "  -4.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "312.02"
$ws.Range("E21").Value = "  -1.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.41"
$ws.Range("E22").Value = "  -1.93%  "
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.52"
$ws.Range("E24").Value = "  -1.44%  "
$ws.Range("E25").Value = "  -1.56%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.03"
$ws.Range("E27").Value = "  -3.71%  "
$ws.Range("E28").Value = "  -5.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "169.78"
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.71"
$ws.Range("E30").Value = "  -4.13%  "
$ws.Range("D31").Value = "0.0₃0720"
$ws.Range("E31").Value = "  -2.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.74"
$ws.Range("E32").Value = "  -2.24%  "
$ws.Range("E33").Value = "  -3.93%  "
$ws.Range("E34").Value = "  -3.99%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.76"
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("E38").Value = "  -4.27%  "
$ws.Range("E39").Value = "  -3.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "38.42"
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("E41").Value = "  -4.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "141.28"
$ws.Range("E42").Value = "  -1.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "288.82"
$ws.Range("E43").Value = "  -5.27%  "
$ws.Range("E44").Value = "  -1.32%  "
$ws.Range("E45").Value = "  -0.26%  "
$ws.Range("E46").Value = "  -1.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.554"
$ws.Range("E47").Value = "  -0.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.11"
$ws.Range("E48").Value = "  -3.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0210"
$ws.Range("E49").Value = "  -2.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.92"
$ws.Range("E50").Value = "  -1.13%  "
$ws.Range("E51").Value = "  -0.65%  "
